$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update row 8 metrics (ano 2025) with refreshed data
$ws.Range("C8").Value = 970
$ws.Range("E8").Value = 807
$ws.Range("G8").Value = 83.1958762886598
$ws.Range("H8").Value = 16.80412371134021
